$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4983071088790894
$ws.Range("B1").Value = 3.439460277557373
$ws.Range("C1").Value = 4.420679092407227
$ws.Range("D1").Value = 2.50184154510498
$ws.Range("E1").Value = 1.062307238578796
